$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCS")

# The sheet has two mini-tables side by side, each with a "Region" column
# that is being removed (region handling moved to a separate "trans" file).
# Left table's Region column is A; right table's Region column is O.
# Delete column A first (everything right of it shifts left by one), then
# the right table's former column O - now at N - is deleted as well.
$ws.Range("A:A").Delete()
$ws.Range("N:N").Delete()

$wb.Save()
